$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[4.51626377967999, 8.83042552615745]"
$ws.Range("M2").Value = 0.000000002648890440326568
$ws.Range("N2").Value = 0.000000005297780880653136
$ws.Range("T2").Value = "[7.9081988174238225, 10.469603729445993]"

# Row 3 updates
$ws.Range("L3").Value = "[4.38478583655799, 9.703372155825267]"
$ws.Range("M3").Value = 0.0000003307476914571339
$ws.Range("N3").Value = 0.0000003307476914571339
$ws.Range("P3").Value = "[0.15723686954903915, 0.9622896416401163]"
$ws.Range("Q3").Value = 0.006561914776339606
$ws.Range("R3").Value = 0.006561914776339606
$ws.Range("T3").Value = "[7.6030449432077365, 10.440736732201016]"
$ws.Range("X3").Value = 19.62990990991007
$ws.Range("Y3").Value = 22.5999199199201
